$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Supplier info block (B3:B5) ---
$ws.Range("B3").Value = "宁波泰丰机械有限公司"
$ws.Range("B4").Value = 13605889085
$ws.Range("B5").Value = "郑小姐"

# --- Line item row 7: description / qty / unit price ---
$ws.Range("C7").Value = "白色针+白色（TPR）胶皮参考denman样品，手柄284U蓝色透明弹性漆，背面印白色logo 包装方式：纸卡包装系扎带，包装背面贴fba不干胶，套opp袋，注意单箱不要超过100个"

$ws.Range("D7").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("D7").Value = 25569.33333333333

$ws.Range("E7").Value = 5.14
